$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GVA")

$ws.Range("B2").Value = 0.1255218910902343
$ws.Range("C2").Value = 0.7467324375882298
$ws.Range("D2").Value = 0.8624827145106858
$ws.Range("E2").Value = 0.9286994748090933
$ws.Range("F2").Value = 0.9468550248872103
$ws.Range("G2").Value = 18
